# Update the "Förändrad" (Changed) date column (C) for all data rows.
# Every row from 2 to 411 currently holds the serial date 45205 (2023-10-06)
# and needs to be bumped by one day to 45206 (2023-10-07).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 411 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
